# Fix Training Data Issue
# The "Date" column (BF) for every data row held the sheet's own file-name-derived
# label ("2-17-2007-08") instead of the actual game date. Correct it to the
# real ISO date "2008-02-17" for all 30 data rows (BF2:BF31).
#
# The literal text is assigned with a leading apostrophe so Excel treats it as
# plain text instead of auto-converting the "YYYY-MM-DD"-looking string into a
# date serial number; the Style reset afterwards drops the quote-prefix
# formatting flag so the cells end up with no extra number formatting applied,
# matching their original (unstyled) appearance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")
$dateRange.Value = "'2008-02-17"
$dateRange.Style = "Normal"
